# "Bug report.xlsx" - fill in row 3 with a real bug report (it used to be
# a verbatim copy of the row-2 placeholder text) and record who/when.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$bullet = [char]0x2022
$bugDescription = "50024: Corner path failure" + "`n" + "Description" + "`n" + `
    "Task: T_ROB1.Corner path executed as stop point due to some of the following reasons:          " + `
    $bullet + " Time delay. " + $bullet + " Closely programmed points. " + $bullet + `
    " System requires high CPU load.Program ref. /SERVER/main/MoveAbsJ/377."

$ws.Range("A3").Value = "InitialPos in Robot_YUMI.py"
$ws.Range("B3").Value = $bugDescription
$ws.Range("E3").Value = "Jiantao Shen"
$ws.Range("D3").Value = "Disconnect the robotic arm Yumi."
$ws.Range("C3").Value = "The robot arm should set all joints to specific angles."

# Date for last update
$ws.Range("F3").NumberFormat = "mm-dd-yy"
$ws.Range("F3").Value = 45198

# The filled-in bug description is much longer than the placeholder text,
# so the row needs to grow to fit it.
$ws.Rows.Item(3).RowHeight = 223.2

# Selection moved onto the author cell we just filled in.
$ws.Range("E3").Select() | Out-Null

$wb.Save()
